$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 33.64264157748815
$ws.Range("G2").Value = 33.04822614340361
$ws.Range("H2").Value = 34.26255084326154
$ws.Range("I2").Value = 3.852848666237324
$ws.Range("J2").Value = 3.821485841668898
$ws.Range("K2").Value = 3.884203944790998
$ws.Range("L2").Value = 0.2696032038637817
$ws.Range("M2").Value = 0.2672772571421074
$ws.Range("N2").Value = 0.271938621608267

# Row 3
$ws.Range("F3").Value = 0.1004515602882827
$ws.Range("G3").Value = 0.01119481366378316
$ws.Range("H3").Value = 0.2120291421887047
$ws.Range("I3").Value = 0.09155061170246738
$ws.Range("J3").Value = 0.01019581812649383
$ws.Range("K3").Value = 0.193241689418152
$ws.Range("L3").Value = 0.1074134381666064
$ws.Range("M3").Value = 0.01202148096992459
$ws.Range("N3").Value = 0.2264050143074009

# Row 4
$ws.Range("F4").Value = 33.74309313777643
$ws.Range("G4").Value = 33.05942095706739
$ws.Range("H4").Value = 34.47457998545025
$ws.Range("I4").Value = 3.944399277939791
$ws.Range("J4").Value = 3.831681659795392
$ws.Range("K4").Value = 4.07744563420915
$ws.Range("L4").Value = 0.377016642030388
$ws.Range("M4").Value = 0.2792987381120319
$ws.Range("N4").Value = 0.4983436359156679
